$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 21
$ws.Range("B2").Value = 35
$ws.Range("B4").Value = 25
$ws.Range("B5").Value = 87
$ws.Range("B6").Value = 45
$ws.Range("B7").Value = 12
$ws.Range("B8").Value = 89
$ws.Range("B10").Value = 90
